$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').NumberFormat = '@'
$ws.Range('D2').Value = '35.845.04'
$ws.Range('E2').Value = '  +3.94%  '
$ws.Range('D3').NumberFormat = '@'
$ws.Range('D3').Value = '1.875.67'
$ws.Range('E3').Value = '  +3.44%  '
$ws.Range('D4').NumberFormat = '@'
$ws.Range('D4').Value = '1.01'
$ws.Range('E4').Value = '  +0.37%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '232.54'
$ws.Range('E5').Value = '  +2.94%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '0.615'
$ws.Range('E6').Value = '  +3.32%  '
$ws.Range('E7').Value = '  +0.23%  '
$ws.Range('D8').NumberFormat = '@'
$ws.Range('D8').Value = '42.29'
$ws.Range('E8').Value = '  +10.17%  '
$ws.Range('D9').NumberFormat = '@'
$ws.Range('D9').Value = '0.312'
$ws.Range('E9').Value = '  +7.73%  '
$ws.Range('D10').NumberFormat = '@'
$ws.Range('D10').Value = '0.0700'
$ws.Range('E10').Value = '  +3.62%  '
$ws.Range('E11').Value = '  +4.05%  '
$ws.Range('D12').NumberFormat = '@'
$ws.Range('D12').Value = '2.144.65'
$ws.Range('E12').Value = '  +3.36%  '
$ws.Range('D13').NumberFormat = '@'
$ws.Range('D13').Value = '11.70'
$ws.Range('E13').Value = '  +4.53%  '
$ws.Range('B14').Value = 'WrappedEther'
$ws.Range('C14').Value = 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth'
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '1.873.98'
$ws.Range('E14').Value = '  +2.96%  '
$ws.Range('B15').Value = 'Polygon'
$ws.Range('C15').Value = 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic'
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.687'
$ws.Range('E15').Value = '  +8.45%  '
$ws.Range('D16').NumberFormat = '@'
$ws.Range('D16').Value = '4.80'
$ws.Range('E16').Value = '  +8.60%  '
$ws.Range('D17').NumberFormat = '@'
$ws.Range('D17').Value = '35.888.56'
$ws.Range('E17').Value = '  +4.13%  '
$ws.Range('D18').NumberFormat = '@'
$ws.Range('D18').Value = '70.81'
$ws.Range('E18').Value = '  +3.63%  '
$ws.Range('E19').Value = '  +4.40%  '
$ws.Range('D20').NumberFormat = '@'
$ws.Range('D20').Value = '248.09'
$ws.Range('E20').Value = '  +2.25%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '12.54'
$ws.Range('E21').Value = '  +11.54%  '
$ws.Range('D22').NumberFormat = '@'
$ws.Range('D22').Value = '4.82'
$ws.Range('E22').Value = '  +16.78%  '
$ws.Range('E23').Value = '  +0.24%  '
$ws.Range('E24').Value = '  +1.34%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '171.45'
$ws.Range('E25').Value = '  +0.72%  '
$ws.Range('E26').Value = '  +3.53%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '18.04'
$ws.Range('D28').NumberFormat = '@'
$ws.Range('D28').Value = '0.123'
$ws.Range('E28').Value = '  +2.17%  '
$ws.Range('D29').NumberFormat = '@'
$ws.Range('D29').Value = '1.45'
$ws.Range('E29').Value = '  +18.13%  '
$ws.Range('D30').NumberFormat = '@'
$ws.Range('D30').Value = '1.01'
$ws.Range('E30').Value = '  +0.43%  '
$ws.Range('D31').NumberFormat = '@'
$ws.Range('D31').Value = '3.333.14'
$ws.Range('E31').Value = '  +37.18%  '
$ws.Range('B32').Value = 'Hedera'
$ws.Range('C32').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '0.0551'
$ws.Range('E32').Value = '  +6.43%  '
$ws.Range('B33').Value = 'Filecoin'
$ws.Range('C33').Value = 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil'
$ws.Range('D33').NumberFormat = '@'
$ws.Range('D33').Value = '3.98'
$ws.Range('E33').Value = '  +4.90%  '
$ws.Range('E34').Value = '  +6.75%  '
$ws.Range('E35').Value = '  +5.37%  '
$ws.Range('D36').NumberFormat = '@'
$ws.Range('D36').Value = '98.69'
$ws.Range('E36').Value = '  +20.86%  '
$ws.Range('E37').Value = '  +7.30%  '
$ws.Range('D38').NumberFormat = '@'
$ws.Range('D38').Value = '2.53'
$ws.Range('E38').Value = '  +7.65%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '1.367.08'
$ws.Range('E39').Value = '  +0.21%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '1.09'
$ws.Range('E40').Value = '  +3.24%  '
$ws.Range('E41').Value = '  +5.91%  '
$ws.Range('E42').Value = '  +8.45%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '15.12'
$ws.Range('E43').Value = '  +9.50%  '
$ws.Range('E44').Value = '  +2.96%  '
$ws.Range('E45').Value = '  +2.24%  '
$ws.Range('E46').Value = '  +1.27%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '6.33'
$ws.Range('E47').Value = '  +9.73%  '
$ws.Range('E48').Value = '  +2.01%  '
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '2.043.05'
$ws.Range('E49').Value = '  +3.42%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '105.47'
$ws.Range('E50').Value = '  +3.35%  '
$ws.Range('E51').Value = '  +0.20%  '
